$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "timestamp" column (O) for every data row (2..398)
# from "2023-01-07 06:49:25" to "2023-01-07 12:54:35"
$ws.Range("O2:O398").Value = "2023-01-07 12:54:35"
